$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add a new row to the "Tabela2" table (extends table ref + dimension)
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# 2. Copy the formatting of the row above (centered style, borders) onto the new row
$ws.Range("C14:E14").Copy()
$ws.Range("C15:E15").PasteSpecial(-4122)

# 3. Fill in the new row's data
$ws.Range("B15").Value = "FUNCIONALIDADE DE CADASTRAR USUÁRIOS"
$ws.Range("C15").Value = 10
$ws.Range("D15").Value = 7

# 4. Extend the conditional formatting on column E down to the new row
$fcs = $ws.Range("E4:E14").FormatConditions
$fcs.Item(1).ModifyAppliesToRange($ws.Range("E4:E15"))

# 5. Extend the data validation list on column E down to the new row
$valRange = $ws.Range("E4:E15")
$valRange.Validation.Delete()
$valRange.Validation.Add(3, 1, 1, "=`$G`$4:`$G`$6")

# 6. Recreate (then remove) the OK/ANDAMENTO/ATRASADO conditional formatting rules
#    on a scratch range, matching the extra style records left behind in the workbook
$scratch = $ws.Range("G24:G26")
$t1 = $scratch.FormatConditions.Add(1, 3, '="OK"')
$t1.Interior.Color = 9364099
$t2 = $scratch.FormatConditions.Add(1, 3, '="ANDAMENTO"')
$t2.Interior.Color = 65535
$t3 = $scratch.FormatConditions.Add(1, 3, '="ATRASADO"')
$t3.Interior.Color = 255
$scratch.FormatConditions.Item(1).Delete()
$scratch.FormatConditions.Item(1).Delete()
$scratch.FormatConditions.Item(1).Delete()

# 7. Leave the selection where the author's session ended up
$ws.Range("G24").Select() | Out-Null
